$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 2.5

$ws.Range("D12").Select() | Out-Null
